$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("B3").Value = "Pedro Francisco-MTRM"

# Row 4
$ws.Range("B4").Value = "Pedro Francisco-MTRM"
$ws.Range("F4").Value = "João Rodrigues-CAD"

# Row 6
$ws.Range("C6").Value = "-"
$ws.Range("D6").Value = "-"
$ws.Range("E6").Value = "Euclides-Mecanica material"
$ws.Range("F6").Value = "João Rodrigues-CAD"

# Row 7
$ws.Range("C7").Value = "-"
$ws.Range("D7").Value = "-"
$ws.Range("E7").Value = "Euclides-Mecanica material"
$ws.Range("F7").Value = "-"

# Row 8
$ws.Range("E8").Value = "-"
